$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull/push of data + mean calculation
$ws.Range("F8").Value = 3
$ws.Range("F12").Value = -5
$ws.Range("F14").Value = -3
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = -13
$ws.Range("F21").Value = -5
